$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.66
$ws.Range("G2").Value = 1.68
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 4.3
$ws.Range("K2").Value = 4.4
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 4.7
$ws.Range("O2").Value = 1.26
$ws.Range("P2").Value = 2.24
$ws.Range("Q2").Value = 1.79
$ws.Range("R2").Value = 1.48
$ws.Range("S2").Value = 2.98
$ws.Range("W2").Value = 2.46
$ws.Range("X2").Value = 20
$ws.Range("Y2").Value = 23
$ws.Range("Z2").Value = 50
$ws.Range("AA2").Value = 150
$ws.Range("AB2").Value = 9.8
$ws.Range("AC2").Value = 9.6
$ws.Range("AI2").Value = 75
$ws.Range("AJ2").Value = 16.5
$ws.Range("AK2").Value = 16
$ws.Range("AM2").Value = 100
$ws.Range("AN2").Value = 8.6
$ws.Range("AO2").Value = 80

$ws.Range("L3").Value = 1.31
$ws.Range("M3").Value = 1.03

$ws.Range("F4").Value = 1.38
$ws.Range("G4").Value = 1.42
$ws.Range("I4").Value = 16.5
$ws.Range("K4").Value = 5
$ws.Range("L4").Value = 1.5
$ws.Range("M4").Value = 1.09
$ws.Range("N4").Value = 3
$ws.Range("O4").Value = 1.43
$ws.Range("Q4").Value = 2.28
$ws.Range("R4").Value = 1.23
$ws.Range("S4").Value = 4.5
$ws.Range("T4").Value = 2.72
$ws.Range("W4").Value = 3.35
$ws.Range("X4").Value = 11.5
$ws.Range("Y4").Value = 29
$ws.Range("AC4").Value = 11.5

$ws.Range("F6").Value = 8.4
$ws.Range("G6").Value = 8.8
$ws.Range("H6").Value = 1.48
$ws.Range("I6").Value = 1.49
$ws.Range("J6").Value = 4.7
$ws.Range("K6").Value = 4.9
$ws.Range("L6").Value = 1.4
$ws.Range("N6").Value = 3.95
$ws.Range("O6").Value = 1.31
$ws.Range("P6").Value = 2.02
$ws.Range("Q6").Value = 1.95
$ws.Range("R6").Value = 1.38
$ws.Range("S6").Value = 3.4
$ws.Range("T6").Value = 2.08
$ws.Range("U6").Value = 1.83
$ws.Range("V6").Value = 3
$ws.Range("X6").Value = 30
$ws.Range("Y6").Value = 7.6
$ws.Range("AB6").Value = 26
$ws.Range("AF6").Value = 70
$ws.Range("AG6").Value = 32
$ws.Range("AI6").Value = 40
$ws.Range("AJ6").Value = 310
$ws.Range("AL6").Value = 140
$ws.Range("AN6").Value = 200
$ws.Range("AO6").Value = 8.4

$ws.Range("F7").Value = 1.06

